$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample-size data after filtration (sample index 1..11, B/C/D counts)
$data = @(
    @(1, 14.0, 11.0, 29.0),
    @(2, 14.0, 11.0, 28.0),
    @(3, 14.0, 9.0, 20.0),
    @(4, 14.0, 11.0, 29.0),
    @(5, 14.0, 11.0, 29.0),
    @(6, 14.0, 11.0, 29.0),
    @(7, 12.0, 5.0, 13.0),
    @(8, 12.0, 3.0, 12.0),
    @(9, 12.0, 11.0, 29.0),
    @(10, 14.0, 11.0, 29.0),
    @(11, 12.0, 7.0, 23.0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = "'" + [string]$vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$wb.Save()
